# Update "想去人数" (interested-count) values in the "展览" and "全部类型" sheets.
# These are simple value refreshes scraped at a later point in time.

$wb = $excel.ActiveWorkbook

# Cell -> new value, keyed by sheet name. Row numbers differ between the two
# sheets because "全部类型" has one extra row inserted part-way through.
$sheetChanges = @{
    "展览" = @{
        "F2" = 200
        "F6" = 18509
        "F7" = 385
        "F8" = 276
        "F9" = 1075
        "F10" = 6935
        "F11" = 192
        "F12" = 695
        "F15" = 120
        "F16" = 77
        "F17" = 222
        "F18" = 169
        "F20" = 279
        "F21" = 60
        "F26" = 286
        "F27" = 1012
        "F28" = 6
        "F29" = 134
        "F30" = 5186
        "F32" = 8
        "F33" = 63
        "F36" = 12182
        "F38" = 15
        "F40" = 217
        "F41" = 298
    }
    "全部类型" = @{
        "F2" = 200
        "F6" = 18509
        "F7" = 385
        "F8" = 276
        "F9" = 1075
        "F10" = 6935
        "F11" = 192
        "F12" = 695
        "F15" = 120
        "F16" = 77
        "F17" = 222
        "F18" = 169
        "F20" = 279
        "F21" = 60
        "F26" = 286
        "F27" = 1012
        "F28" = 6
        "F29" = 134
        "F30" = 5186
        "F33" = 8
        "F35" = 63
        "F38" = 12182
        "F40" = 15
        "F42" = 217
        "F43" = 298
    }
}

foreach ($sheetName in $sheetChanges.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetChanges[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
